$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 19 (Leve Item ID 7015)
$ws.Range("H19").Value = 1082.1786
$ws.Range("I19").Value = 613.7
$ws.Range("J19").Value = 1342.4445
$ws.Range("K19").Value = 613.7
$ws.Range("L19").Value = 1342.4445
$ws.Range("M19").Value = -438.7
$ws.Range("N19").Value = -1692.4445

# Row 21 (Leve Item ID 2149)
$ws.Range("H21").Value = 4699.3335
$ws.Range("I21").Value = 4699.3335
$ws.Range("K21").Value = 4699.3335
$ws.Range("M21").Value = -4231.3335

# Row 23 (Leve Item ID 2149)
$ws.Range("H23").Value = 4699.3335
$ws.Range("I23").Value = 4699.3335
$ws.Range("K23").Value = 4699.3335
$ws.Range("M23").Value = -4465.3335

# Row 31 (Leve Item ID 4576)
$ws.Range("H31").Value = 43403.332
$ws.Range("I31").Value = 43403.332
$ws.Range("K31").Value = 130209.996
$ws.Range("M31").Value = -129979.996

# Row 32 (Leve Item ID 5484)
$ws.Range("H32").Value = 6333.3335
$ws.Range("I32").Value = 0
$ws.Range("K32").Value = 0
$ws.Range("M32").ClearContents()

# Row 41 (Leve Item ID 5478)
$ws.Range("H41").Value = 284.75
$ws.Range("I41").Value = 284.75
$ws.Range("K41").Value = 284.75
$ws.Range("M41").Value = 155.25

# Row 44 (Leve Item ID 1971)
$ws.Range("H44").Value = 10000
$ws.Range("J44").Value = 10000
$ws.Range("L44").Value = 10000
$ws.Range("N44").Value = -10924

# Row 70 (Leve Item ID 12604)
$ws.Range("H70").Value = 5628.4287
$ws.Range("I70").Value = 6106.6665
$ws.Range("J70").Value = 5076.615
$ws.Range("K70").Value = 18319.9995
$ws.Range("L70").Value = 15229.845
$ws.Range("M70").Value = -18049.9995
$ws.Range("N70").Value = -15769.845

# Row 73 (Leve Item ID 12604)
$ws.Range("H73").Value = 5628.4287
$ws.Range("I73").Value = 6106.6665
$ws.Range("J73").Value = 5076.615
$ws.Range("K73").Value = 18319.9995
$ws.Range("L73").Value = 15229.845
$ws.Range("M73").Value = -17383.9995
$ws.Range("N73").Value = -17101.845

# Row 86 (Leve Item ID 12603)
$ws.Range("H86").Value = 4835.7617
$ws.Range("I86").Value = 5860.6
$ws.Range("J86").Value = 3904.0908
$ws.Range("K86").Value = 5860.6
$ws.Range("L86").Value = 3904.0908
$ws.Range("M86").Value = -4737.6
$ws.Range("N86").Value = -6150.0908

# Row 89 (Leve Item ID 12603)
$ws.Range("H89").Value = 4835.7617
$ws.Range("I89").Value = 5860.6
$ws.Range("J89").Value = 3904.0908
$ws.Range("K89").Value = 29303
$ws.Range("L89").Value = 19520.454
$ws.Range("M89").Value = -23687
$ws.Range("N89").Value = -30752.454

# Row 98 (Leve Item ID 36237)
$ws.Range("H98").Value = 942.5
$ws.Range("I98").Value = 942.5
$ws.Range("K98").Value = 942.5
$ws.Range("M98").Value = 555.5

# Row 122 (Leve Item ID 36237)
$ws.Range("H122").Value = 942.5
$ws.Range("I122").Value = 942.5
$ws.Range("K122").Value = 2827.5
$ws.Range("M122").Value = -377.5

# Row 132 (Leve Item ID 44049)
$ws.Range("H132").Value = 5662
$ws.Range("I132").Value = 1433.3334
$ws.Range("J132").Value = 8199.200000000001
$ws.Range("K132").Value = 4300.0002
$ws.Range("L132").Value = 24597.6
$ws.Range("M132").Value = -1770.0002
$ws.Range("N132").Value = -29657.6

$ws = $wb.Worksheets.Item("ARM")
# Row 2 (Leve Item ID 27713)
$ws.Range("H2").Value = 41667796
$ws.Range("I2").Value = 41667796
$ws.Range("J2").Value = 0
$ws.Range("K2").Value = 41667796
$ws.Range("L2").Value = 0
$ws.Range("N2").Value = -41667683
$ws.Range("M2").ClearContents()

# Row 32 (Leve Item ID 44147)
$ws.Range("H32").Value = 239343.53
$ws.Range("I32").Value = 781.8919
$ws.Range("K32").Value = 781.8919
$ws.Range("M32").Value = -494.8919

# Row 46 (Leve Item ID 3498)
$ws.Range("H46").Value = 16466.666
$ws.Range("J46").Value = 9950
$ws.Range("L46").Value = 9950
$ws.Range("N46").Value = -10588

# Row 116 (Leve Item ID 27713)
$ws.Range("H116").Value = 41667796
$ws.Range("I116").Value = 41667796
$ws.Range("J116").Value = 0
$ws.Range("K116").Value = 41667796
$ws.Range("L116").Value = 0
$ws.Range("N116").Value = -41665502
$ws.Range("M116").ClearContents()

$ws = $wb.Worksheets.Item("BSM")
# Row 3 (Leve Item ID 27713)
$ws.Range("H3").Value = 41667796
$ws.Range("I3").Value = 41667796
$ws.Range("J3").Value = 0
$ws.Range("K3").Value = 41667796
$ws.Range("L3").Value = 0
$ws.Range("N3").Value = -41667682
$ws.Range("M3").ClearContents()

$ws = $wb.Worksheets.Item("CRP")
# Row 22 (Leve Item ID 5367)
$ws.Range("H22").Value = 2977.182
$ws.Range("I22").Value = 3154.9
$ws.Range("J22").Value = 1200
$ws.Range("K22").Value = 3154.9
$ws.Range("L22").Value = 1200
$ws.Range("M22").Value = -2804.9
$ws.Range("N22").Value = -1900

# Row 31 (Leve Item ID 44023)
$ws.Range("H31").Value = 4410.2
$ws.Range("I31").Value = 4031
$ws.Range("J31").Value = 4663
$ws.Range("K31").Value = 4031
$ws.Range("L31").Value = 4663
$ws.Range("M31").Value = -3736
$ws.Range("N31").Value = -5253

# Row 34 (Leve Item ID 44023)
$ws.Range("H34").Value = 4410.2
$ws.Range("I34").Value = 4031
$ws.Range("J34").Value = 4663
$ws.Range("K34").Value = 4031
$ws.Range("L34").Value = 4663
$ws.Range("M34").Value = -3829
$ws.Range("N34").Value = -5067

# Row 132 (Leve Item ID 44019)
$ws.Range("H132").Value = 1585.579
$ws.Range("I132").Value = 1389.7646
$ws.Range("K132").Value = 4169.293799999999
$ws.Range("M132").Value = -1639.293799999999

$ws = $wb.Worksheets.Item("CUL")
# Row 38 (Leve Item ID 4860)
$ws.Range("H38").Value = 259.4
$ws.Range("I38").Value = 35
$ws.Range("J38").Value = 315.5
$ws.Range("K38").Value = 105
$ws.Range("L38").Value = 946.5
$ws.Range("M38").Value = 242
$ws.Range("N38").Value = -1640.5

# Row 55 (Leve Item ID 4733)
$ws.Range("H55").Value = 3578.6
$ws.Range("J55").Value = 3974.5
$ws.Range("L55").Value = 11923.5
$ws.Range("N55").Value = -12277.5

# Row 109 (Leve Item ID 27854)
$ws.Range("H109").Value = 953.25
$ws.Range("I109").Value = 661
$ws.Range("K109").Value = 1983
$ws.Range("M109").Value = -943

$ws = $wb.Worksheets.Item("GSM")
# Row 2 (Leve Item ID 5062)
$ws.Range("H2").Value = 262.125
$ws.Range("I2").Value = 57.714287
$ws.Range("J2").Value = 421.1111
$ws.Range("K2").Value = 57.714287
$ws.Range("L2").Value = 421.1111
$ws.Range("M2").Value = 55.285713
$ws.Range("N2").Value = -647.1111000000001

# Row 80 (Leve Item ID 12521)
$ws.Range("H80").Value = 1664.75
$ws.Range("J80").Value = 1554.5
$ws.Range("L80").Value = 1554.5
$ws.Range("N80").Value = -3550.5

# Row 83 (Leve Item ID 12521)
$ws.Range("H83").Value = 1664.75
$ws.Range("J83").Value = 1554.5
$ws.Range("L83").Value = 7772.5
$ws.Range("N83").Value = -17756.5

$ws = $wb.Worksheets.Item("LTW")
# Row 22 (Leve Item ID 5277)
$ws.Range("H22").Value = 1547
$ws.Range("J22").Value = 1729.3334
$ws.Range("L22").Value = 1729.3334
$ws.Range("N22").Value = -2319.3334

# Row 27 (Leve Item ID 5277)
$ws.Range("H27").Value = 1547
$ws.Range("J27").Value = 1729.3334
$ws.Range("L27").Value = 1729.3334
$ws.Range("N27").Value = -1943.3334

$ws = $wb.Worksheets.Item("WVR")
# Row 31 (Leve Item ID 3052)
$ws.Range("H31").Value = 18316.666
$ws.Range("I31").Value = 14850
$ws.Range("K31").Value = 14850
$ws.Range("M31").Value = -14502

# Row 48 (Leve Item ID 3140)
$ws.Range("H48").Value = 42500
$ws.Range("J48").Value = 42500
$ws.Range("L48").Value = 42500
$ws.Range("N48").Value = -43638

# Row 132 (Leve Item ID 44029)
$ws.Range("H132").Value = 6357.067
$ws.Range("I132").Value = 3200.2222
$ws.Range("J132").Value = 11092.333
$ws.Range("K132").Value = 9600.6666
$ws.Range("L132").Value = 33276.999
$ws.Range("M132").Value = -7070.6666
$ws.Range("N132").Value = -38336.999
